$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Row 13: the previously blank cells B13:K13 and P13 now hold "nan"
$ws.Range("B13:K13").Value = "nan"
$ws.Range("P13").Value = "nan"

# Row 14: new service-log entry for Card24
# A14 holds the card number as text (matches A2:A13 which are all text "24")
$ws.Range("A14").Value = "'24"

# L14/N14/O14 carry the new event's date / correction / serviced-by text
$ws.Range("L14").Value = "23\9\2024"
$ws.Range("N14").Value = "تم عمل صيانه ربع سنويه"
$ws.Range("O14").Value = "تيم العمل"

# The remaining row-14 cells (B:K, M, P) stay blank but must still exist
# as real (empty) cells so the sheet dimension extends to row 14.
$blankCols = @("B","C","D","E","F","G","H","I","J","K","M","P")
foreach ($col in $blankCols) {
    $ws.Range($col + "14").Style = "Normal"
}
